$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 256.86365
$ws.Range("I33").Value = 173.85715
$ws.Range("J33").Value = 2000
$ws.Range("K33").Value = 173.85715
$ws.Range("L33").Value = 2000
$ws.Range("M33").Value = 55.14285000000001
$ws.Range("N33").Value = -2458

$ws.Range("H64").Value = 6957.8237
$ws.Range("I64").Value = 6964.778
$ws.Range("J64").Value = 6950
$ws.Range("K64").Value = 6964.778
$ws.Range("L64").Value = 6950
$ws.Range("M64").Value = -6716.778
$ws.Range("N64").Value = -7446

$ws.Range("H67").Value = 6957.8237
$ws.Range("I67").Value = 6964.778
$ws.Range("J67").Value = 6950
$ws.Range("K67").Value = 6964.778
$ws.Range("L67").Value = 6950
$ws.Range("M67").Value = -6106.778
$ws.Range("N67").Value = -8666

$ws.Range("H74").Value = 5858.6
$ws.Range("I74").Value = 3976.7273
$ws.Range("K74").Value = 3976.7273
$ws.Range("M74").Value = -3040.7273

$ws.Range("H77").Value = 5858.6
$ws.Range("I77").Value = 3976.7273
$ws.Range("K77").Value = 19883.6365
$ws.Range("M77").Value = -15203.6365

$ws.Range("H86").Value = 2832.875
$ws.Range("I86").Value = 2202.2
$ws.Range("K86").Value = 2202.2
$ws.Range("M86").Value = -1079.2

$ws.Range("H89").Value = 2832.875
$ws.Range("I89").Value = 2202.2
$ws.Range("K89").Value = 11011
$ws.Range("M89").Value = -5395

$ws.Range("H132").Value = 67362.61
$ws.Range("I132").Value = 74458
$ws.Range("K132").Value = 223374
$ws.Range("M132").Value = -220844

$ws.Range("H138").Value = 2270.244
$ws.Range("I138").Value = 1700.238
$ws.Range("J138").Value = 2868.75
$ws.Range("K138").Value = 5100.714
$ws.Range("L138").Value = 8606.25
$ws.Range("M138").Value = 39.28600000000006
$ws.Range("N138").Value = -18886.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 2013
$ws.Range("J15").Value = 2013
$ws.Range("L15").Value = 2013
$ws.Range("N15").Value = -2713

$ws.Range("H39").Value = 14999
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H61").Value = 1391339.9
$ws.Range("I61").Value = 1589421.1
$ws.Range("K61").Value = 1589421.1
$ws.Range("M61").Value = -1589209.1

$ws.Range("H74").Value = 2062.068
$ws.Range("I74").Value = 745.3333
$ws.Range("J74").Value = 7987.375
$ws.Range("K74").Value = 745.3333
$ws.Range("L74").Value = 7987.375
$ws.Range("M74").Value = 128.6667
$ws.Range("N74").Value = -9735.375

$ws.Range("H77").Value = 2062.068
$ws.Range("I77").Value = 745.3333
$ws.Range("J77").Value = 7987.375
$ws.Range("K77").Value = 3726.6665
$ws.Range("L77").Value = 39936.875
$ws.Range("M77").Value = 641.3334999999997
$ws.Range("N77").Value = -48672.875

$ws.Range("H102").Value = 21283.723
$ws.Range("I102").Value = 23716.125
$ws.Range("K102").Value = 23716.125
$ws.Range("M102").Value = -22094.125

$ws.Range("H122").Value = 3921.8635
$ws.Range("I122").Value = 3945
$ws.Range("K122").Value = 11835
$ws.Range("M122").Value = -9385

$ws.Range("H124").Value = 29999.5
$ws.Range("I124").Value = 29999
$ws.Range("K124").Value = 29999
$ws.Range("M124").Value = -25089

$ws.Range("H125").Value = 70178.5
$ws.Range("I125").Value = 64999
$ws.Range("J125").Value = 71905
$ws.Range("K125").Value = 64999
$ws.Range("L125").Value = 71905
$ws.Range("M125").Value = -60079
$ws.Range("N125").Value = -81745

$ws.Range("H136").Value = 1391339.9
$ws.Range("I136").Value = 1589421.1
$ws.Range("K136").Value = 4768263.300000001
$ws.Range("M136").Value = -4765713.300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 13993
$ws.Range("J31").Value = 13993
$ws.Range("L31").Value = 13993
$ws.Range("N31").Value = -14497

$ws.Range("H33").Value = 24000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 24000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 24000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -24672

$ws.Range("H105").Value = 1999.8
$ws.Range("I105").Value = 1499.75
$ws.Range("K105").Value = 1499.75
$ws.Range("M105").Value = 247.25

$ws.Range("H134").Value = 687439.6
$ws.Range("I134").Value = 930885.25
$ws.Range("J134").Value = 291840.5
$ws.Range("K134").Value = 2792655.75
$ws.Range("L134").Value = 875521.5
$ws.Range("M134").Value = -2790120.75
$ws.Range("N134").Value = -880591.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 30049.5
$ws.Range("J51").Value = 30099
$ws.Range("L51").Value = 30099
$ws.Range("N51").Value = -31571

$ws.Range("H61").Value = 30049.5
$ws.Range("J61").Value = 30099
$ws.Range("L61").Value = 30099
$ws.Range("N61").Value = -30795

$ws.Range("H86").Value = 6750.0527
$ws.Range("I86").Value = 5885
$ws.Range("J86").Value = 8624.333000000001
$ws.Range("K86").Value = 5885
$ws.Range("L86").Value = 8624.333000000001
$ws.Range("M86").Value = -4762
$ws.Range("N86").Value = -10870.333

$ws.Range("H89").Value = 6750.0527
$ws.Range("I89").Value = 5885
$ws.Range("J89").Value = 8624.333000000001
$ws.Range("K89").Value = 29425
$ws.Range("L89").Value = 43121.665
$ws.Range("M89").Value = -23809
$ws.Range("N89").Value = -54353.665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2142.6191
$ws.Range("I114").Value = 221.8
$ws.Range("J114").Value = 2742.875
$ws.Range("K114").Value = 665.4000000000001
$ws.Range("L114").Value = 8228.625
$ws.Range("M114").Value = 2588.6
$ws.Range("N114").Value = -14736.625

$ws.Range("H132").Value = 1691.5
$ws.Range("I132").Value = 984.8182
$ws.Range("K132").Value = 8863.363800000001
$ws.Range("M132").Value = -6333.363800000001

$ws.Range("H140").Value = 4342.75
$ws.Range("J140").Value = 4000
$ws.Range("L140").Value = 12000
$ws.Range("N140").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 789.2353000000001
$ws.Range("I2").Value = 910.0833
$ws.Range("K2").Value = 910.0833
$ws.Range("M2").Value = -797.0833

$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10586

$ws.Range("H38").Value = 24024
$ws.Range("J38").Value = 24024
$ws.Range("L38").Value = 24024
$ws.Range("N38").Value = -24950

$ws.Range("H55").Value = 20944.6
$ws.Range("J55").Value = 29911.5
$ws.Range("L55").Value = 29911.5
$ws.Range("N55").Value = -30565.5

$ws.Range("H63").Value = 27000
$ws.Range("J63").Value = 27000
$ws.Range("L63").Value = 27000
$ws.Range("N63").Value = -28372

$ws.Range("H66").Value = 27000
$ws.Range("J66").Value = 27000
$ws.Range("L66").Value = 81000
$ws.Range("N66").Value = -87864

$ws.Range("H70").Value = 7439.6
$ws.Range("I70").Value = 7439.6
$ws.Range("K70").Value = 7439.6
$ws.Range("M70").Value = -7169.6

$ws.Range("H73").Value = 7439.6
$ws.Range("I73").Value = 7439.6
$ws.Range("K73").Value = 7439.6
$ws.Range("M73").Value = -6503.6

$ws.Range("H80").Value = 272016.4
$ws.Range("I80").Value = 343728.2
$ws.Range("J80").Value = 3097.25
$ws.Range("K80").Value = 343728.2
$ws.Range("L80").Value = 3097.25
$ws.Range("M80").Value = -342730.2
$ws.Range("N80").Value = -5093.25

$ws.Range("H83").Value = 272016.4
$ws.Range("I83").Value = 343728.2
$ws.Range("J83").Value = 3097.25
$ws.Range("K83").Value = 1718641
$ws.Range("L83").Value = 15486.25
$ws.Range("M83").Value = -1713649
$ws.Range("N83").Value = -25470.25

$ws.Range("H102").Value = 3759.3157
$ws.Range("I102").Value = 3260.4707
$ws.Range("K102").Value = 3260.4707
$ws.Range("M102").Value = -1638.4707

$ws.Range("H126").Value = 836265.3
$ws.Range("I126").Value = 1391311.8
$ws.Range("K126").Value = 4173935.4
$ws.Range("M126").Value = -4171465.4

$ws.Range("H132").Value = 17455622
$ws.Range("I132").Value = 25956394
$ws.Range("J132").Value = 6666.737
$ws.Range("K132").Value = 77869182
$ws.Range("L132").Value = 20000.211
$ws.Range("M132").Value = -77866652
$ws.Range("N132").Value = -25060.211

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 348.7
$ws.Range("I9").Value = 185.875
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 185.875
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 38.125
$ws.Range("N9").Value = -1448

$ws.Range("H33").Value = 10978.75
$ws.Range("I33").Value = 10978.75
$ws.Range("K33").Value = 10978.75
$ws.Range("M33").Value = -10688.75

$ws.Range("H46").Value = 1648
$ws.Range("I46").Value = 830.5714
$ws.Range("J46").Value = 2124.8333
$ws.Range("K46").Value = 830.5714
$ws.Range("L46").Value = 2124.8333
$ws.Range("M46").Value = -642.5714
$ws.Range("N46").Value = -2500.8333

$ws.Range("H127").Value = 163428.9
$ws.Range("J127").Value = 163428.9
$ws.Range("L127").Value = 163428.9
$ws.Range("N127").Value = -173348.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 5050
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 100
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 140
$ws.Range("N20").Value = -10480
